# Expanded all rooms' MTRX Macro arrays
# Adds a new relay block (rows 74-75) for the "513A (BIOCOMM RACKS)" room,
# PRO3 device at CP 5A-1 (local), relative relay numbers 49/50, relay
# numbers 1/2, UP/DOWN operation, screen room 567, device PROJ SCREEN.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the DEVICE column first (matches the new-shared-string registration
# order seen in the authored workbook: "PROJ SCREEN" ends up as the first
# newly-added shared string).
$ws.Cells.Item(74, 11).Value = "PROJ SCREEN"            # K74 DEVICE
$ws.Cells.Item(75, 11).Value = "PROJ SCREEN"            # K75 DEVICE

# Row 74 - first row of the new block includes the room/device info columns
$ws.Cells.Item(74, 1).Value = "513A (BIOCOMM RACKS)"   # A74 ROOM NAME
$ws.Cells.Item(74, 2).Value = "CRESTRON"                # B74 MANUFACTURER
$ws.Cells.Item(74, 3).Value = "PRO3"                    # C74 MODEL
$ws.Cells.Item(74, 4).Value = "CP 5A-1"                 # D74 DESIGN DEVICE ID
$ws.Cells.Item(74, 5).Value = "local"                   # E74 IP ADDRESS
$ws.Cells.Item(74, 7).Value = 49                        # G74 RELATIVE RELAY NUMBERS
$ws.Cells.Item(74, 8).Value = 1                         # H74 RELAY NUMBER
$ws.Cells.Item(74, 9).Value = "UP"                      # I74 OPERATION
$ws.Cells.Item(74, 10).Value = 567                      # J74 SCREEN ROOM

# Row 75 - second row of the new block (relay 2 / DOWN)
$ws.Cells.Item(75, 7).Value = 50                        # G75 RELATIVE RELAY NUMBERS
$ws.Cells.Item(75, 8).Value = 2                         # H75 RELAY NUMBER
$ws.Cells.Item(75, 9).Value = "DOWN"                    # I75 OPERATION
$ws.Cells.Item(75, 10).Value = 567                      # J75 SCREEN ROOM

# Update the selection to reflect the last-edited cell, as in the source file
$ws.Range("F74").Select()
